# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The workbook's "K" column (column G) holds per-row computed values that were
# regenerated upstream (recalculated statistic, formerly "Strike#"). This
# writes the freshly-computed K values back into column G for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 2
    6  = 1
    7  = 1
    8  = 3
    9  = 2
    10 = 9
    12 = 0
    13 = 0
    14 = 4
    15 = 1
    16 = 2
    17 = 2
    18 = 0
    19 = 1
    20 = 5
    21 = 0
    22 = 0
    23 = 1
    24 = 2
    25 = 4
    26 = 5
    27 = 2
    28 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
